$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Stash the existing header-row formatting (bold font, grey fill, thin
#    border, centered + wrapped alignment) onto a scratch cell far outside
#    the used range, so it can be restored later with a single formats-only
#    paste (this avoids Excel creating a bunch of throw-away cell styles
#    when re-applying the formatting property by property).
# ---------------------------------------------------------------------------
$scratch = $ws.Range("A100")
$ws.Range("A1").Copy() | Out-Null
$scratch.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Rename header row (A1:U1): "_old" suffix -> "_FV2310", "_new" suffix ->
#    "_FV2404" (the "diff" column header is unchanged).
# ---------------------------------------------------------------------------
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}

# ---------------------------------------------------------------------------
# 3. Clear the header formatting before creating the table: the table
#    creation below would otherwise capture the current header formatting
#    as a "header row" differential style (headerRowDxfId) on the table.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$headerRange.ClearFormats()

# ---------------------------------------------------------------------------
# 4. Convert the data range into an Excel Table (ListObject) named "Table1".
# ---------------------------------------------------------------------------
$range = $ws.Range("A1:U77")
$listObject = $ws.ListObjects.Add(1, $range, [System.Reflection.Missing]::Value, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# ---------------------------------------------------------------------------
# 5. Restore the original header-row formatting from the scratch cell with a
#    single formats-only paste, then clean up the scratch cell.
# ---------------------------------------------------------------------------
$scratch.Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$scratch.Clear()

# ---------------------------------------------------------------------------
# 6. Freeze the header row (pane split after row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$window = $excel.ActiveWindow
$window.FreezePanes = $true
